# "stops resized to picbox" - add newly surveyed bus stops (rows 25-36) to
# the List1 stop table, and normalise the stop-name column's alignment to
# match the rest of the (centered) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# --- new stop rows (X / Y pixel coords + stop name) -----------------------
$ws.Range("B25").Value = 5253
$ws.Range("C25").Value = 4298
$ws.Range("D25").Value = "Uherské Hradiště,Sady"

$ws.Range("B26").Value = 5691
$ws.Range("C26").Value = 4452
$ws.Range("D26").Value = "Uherské Hradiště,Sady,Za Kovárnou"

$ws.Range("B27").Value = 6575
$ws.Range("C27").Value = 4799
$ws.Range("D27").Value = "Uherské Hradiště,Vésky,hor.konec"

$ws.Range("B28").Value = 6778
$ws.Range("C28").Value = 5215
$ws.Range("D28").Value = "Uherské Hradiště,Vésky,střed"

$ws.Range("B29").Value = 7097
$ws.Range("C29").Value = 5875
$ws.Range("D29").Value = "Uherské Hradiště,Míkovice,Hlavní"

$ws.Range("B30").Value = 6172
$ws.Range("C30").Value = 5809
$ws.Range("D30").Value = "Kunovice,,koupaliště"

$ws.Range("B31").Value = 4703
$ws.Range("C31").Value = 5550
$ws.Range("D31").Value = "Kunovice,,Lidická"

$ws.Range("B32").Value = 4362
$ws.Range("C32").Value = 5260
$ws.Range("D32").Value = "Kunovice,,Na Rynku"

$ws.Range("B33").Value = 4334
$ws.Range("C33").Value = 4439
$ws.Range("D33").Value = "Kunovice,,rozc.k žel.st.0.5"

$ws.Range("B34").Value = 4041
$ws.Range("C34").Value = 5412
$ws.Range("D34").Value = "Kunovice,,Panská"

$ws.Range("B35").Value = 3488
$ws.Range("C35").Value = 5181
$ws.Range("D35").Value = "Kunovice,,Na Bělince"

$ws.Range("B36").Value = 2137
$ws.Range("C36").Value = 4915
$ws.Range("D36").Value = "Kunovice,,Let"

# --- normalise stop-name column (D) to center/middle alignment ------------
$dcol = $ws.Range("D1:D43")
$dcol.HorizontalAlignment = -4108
$dcol.VerticalAlignment = -4108

# --- move the saved selection cursor to B37, as in the edited workbook ----
$ws.Range("B37").Select()
